# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- OFF sheet (Home row, A2="H") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 331
$wsOff.Range("C2").Value = 219
$wsOff.Range("D2").Value = 76
$wsOff.Range("E2").Value = 37
$wsOff.Range("G2").Value = 5

# --- DEF sheet (Home row, A2="H") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 404
$wsDef.Range("C2").Value = 272
$wsDef.Range("F2").Value = 7
